# announcements and small slide fixes
# Adds a new "Monday, January 27" announcements slide at the end of the
# deck (slide 4), cloned from the last existing slide (slide 3) so it
# picks up the same layout / placeholder set (title, content, footer).

$p = $ppt.ActivePresentation

# Clone the last slide so we inherit its title/content/footer placeholders
# (incl. the slide-number field) instead of building a bare layout slide.
$s4 = $p.Slides.Item(3).Duplicate()

# ---- Title -----------------------------------------------------------
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Monday, January 27"

# ---- Body content ------------------------------------------------------
$body = $s4.Shapes.Item(2).TextFrame.TextRange

$body.Text = "First assignment is due Tuesday`rDoes anyone still not have SDL2 running?`rGroups: Anybody not have one?`rHW2 is due this Friday (will go over it today)`rLab this week, come prepared to:`rLab this week is just working on homework and getting everyone up to speed.`rNext week will be our first big push on the projects`rToday we will finish sprites quickly and do display trees."

# Paragraph 1: "First " + "assignment" (bold italic) + " is due Tuesday"
$para1 = $body.Paragraphs(1, 1)
$para1.Font.Size = 26
$run = $para1.Characters(7, 10)
$run.Font.Bold = -1
$run.Font.Italic = -1

# Paragraph 2: second level bullet
$para2 = $body.Paragraphs(2, 1)
$para2.Font.Size = 22
$para2.IndentLevel = 2

# Paragraph 3: "Groups" (bold italic) + ": Anybody not have one?"
$para3 = $body.Paragraphs(3, 1)
$para3.Font.Size = 26
$run = $para3.Characters(1, 6)
$run.Font.Bold = -1
$run.Font.Italic = -1

# Paragraph 4: "HW2" (bold italic) + " is due this Friday (will go over it today)"
$para4 = $body.Paragraphs(4, 1)
$para4.Font.Size = 26
$run = $para4.Characters(1, 3)
$run.Font.Bold = -1
$run.Font.Italic = -1

# Paragraph 5: "Lab this week" (bold italic) + ", come prepared to:"
$para5 = $body.Paragraphs(5, 1)
$para5.Font.Size = 26
$run = $para5.Characters(1, 13)
$run.Font.Bold = -1
$run.Font.Italic = -1

# Paragraph 6: second level bullet
$para6 = $body.Paragraphs(6, 1)
$para6.Font.Size = 22
$para6.IndentLevel = 2

# Paragraph 7: second level bullet
$para7 = $body.Paragraphs(7, 1)
$para7.Font.Size = 22
$para7.IndentLevel = 2

# Paragraph 8: plain closing line
$para8 = $body.Paragraphs(8, 1)
$para8.Font.Size = 26
